$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row with the latest scraped data.
$ws.Range("D2").Value = '24.874.96'
$ws.Range("E2").Value = '  +1.33%  '

$ws.Range("D3").Value = '1.690.57'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = '  +0.73%  '

$ws.Range("D5").Value = "'316.33"
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = '  +0.43%  '

$ws.Range("D7").Value = "'0.3950"
$ws.Range("E7").Value = '  +0.86%  '

$ws.Range("D8").Value = "'0.3985"
$ws.Range("E8").Value = '  -1.54%  '

$ws.Range("D9").Value = "'52.75"
$ws.Range("E9").Value = '  -1.67%  '

$ws.Range("D10").Value = "'1.441"
$ws.Range("E10").Value = '  -2.90%  '

$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("D12").Value = "'0.08729"
$ws.Range("E12").Value = '  -0.86%  '

$ws.Range("D13").Value = "'25.51"
$ws.Range("E13").Value = '  -2.41%  '

$ws.Range("D14").Value = "'7.389"
$ws.Range("E14").Value = '  -0.86%  '

$ws.Range("D15").Value = '1.944.95'
$ws.Range("E15").Value = '  +14.56%  '

$ws.Range("D16").Value = "'0.00001339"
$ws.Range("E16").Value = '  -1.45%  '

$ws.Range("D17").Value = "'7.868"
$ws.Range("E17").Value = '  -2.84%  '

$ws.Range("D18").Value = "'94.78"
$ws.Range("E18").Value = '  -2.86%  '

$ws.Range("D19").Value = "'0.07252"
$ws.Range("E19").Value = '  +1.09%  '

$ws.Range("D20").Value = "'20.43"
$ws.Range("E20").Value = '  -0.30%  '

$ws.Range("D21").Value = "'7.177"
$ws.Range("E21").Value = '  -1.70%  '

$ws.Range("D22").Value = "'1.006"
$ws.Range("E22").Value = '  +0.44%  '

$ws.Range("D23").Value = "'14.18"
$ws.Range("E23").Value = '  -1.06%  '

$ws.Range("D24").Value = '24.861.73'
$ws.Range("E24").Value = '  +1.32%  '

$ws.Range("D25").Value = "'2.407"
$ws.Range("E25").Value = '  +3.27%  '

$ws.Range("D26").Value = "'2.837"
$ws.Range("E26").Value = '  -5.88%  '

$ws.Range("D27").Value = "'23.13"
$ws.Range("E27").Value = '  +0.88%  '

$ws.Range("D28").Value = "'6.047"
$ws.Range("E28").Value = '  +3.34%  '

$ws.Range("D29").Value = "'161.93"
$ws.Range("E29").Value = '  -3.68%  '

$ws.Range("D30").Value = "'148.74"
$ws.Range("E30").Value = '  +3.10%  '

$ws.Range("D31").Value = "'8.086"
$ws.Range("E31").Value = '  -3.48%  '

$ws.Range("D32").Value = "'2.604"
$ws.Range("E32").Value = '  +19.88%  '

$ws.Range("D33").Value = '1.941.57'
$ws.Range("E33").Value = '  +3.14%  '

$ws.Range("D34").Value = "'0.08508"
$ws.Range("E34").Value = '  -3.13%  '

$ws.Range("D35").Value = "'0.03110"
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("D36").Value = "'1.036"
$ws.Range("E36").Value = '  -2.07%  '

$ws.Range("D37").Value = "'7.070"
$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("D38").Value = "'0.2866"
$ws.Range("E38").Value = '  +2.18%  '

$ws.Range("D39").Value = "'0.09660"
$ws.Range("E39").Value = '  +5.27%  '

$ws.Range("D40").Value = "'10.78"
$ws.Range("E40").Value = '  -1.15%  '

$ws.Range("D41").Value = "'0.8086"
$ws.Range("E41").Value = '  -7.67%  '

$ws.Range("D42").Value = "'13.92"
$ws.Range("E42").Value = '  -2.07%  '

$ws.Range("D43").Value = "'1.471"
$ws.Range("E43").Value = '  -0.47%  '

$ws.Range("D44").Value = "'16.89"
$ws.Range("E44").Value = '  -2.86%  '

$ws.Range("D45").Value = "'2.629"
$ws.Range("E45").Value = '  -1.42%  '

$ws.Range("D46").Value = "'0.7273"
$ws.Range("E46").Value = '  -3.82%  '

$ws.Range("D47").Value = "'4.222"
$ws.Range("E47").Value = '  -0.88%  '

$ws.Range("D48").Value = "'0.08954"
$ws.Range("E48").Value = '  +9.02%  '

$ws.Range("D49").Value = "'1.382"
$ws.Range("E49").Value = '  -0.72%  '

$ws.Range("D50").Value = "'1.003"
$ws.Range("E50").Value = '  +0.26%  '

$ws.Range("E51").Value = '  -0.41%  '

# The quote-prefix trick above makes Excel apply a "Text" quoted-prefix
# cell style; reapply the default "Normal" style to each affected cell so
# the cell keeps its original (unstyled) appearance.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
